$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the former "Espresso, Decaf" row (previously row 4)
$ws.Range("A2").Formula = '="10029"'
$ws.Range("B2").Formula = '="Equal Exchange - Espresso, Decaf"'
$ws.Range("C2").Formula = '="1"'
$ws.Range("D2").Formula = '="94.00"'
$ws.Range("E2").Formula = '="94.00"'

# Row 3 becomes the former "Cold Brew" row (previously row 5)
$ws.Range("A3").Formula = '="10403"'
$ws.Range("B3").Formula = '="Equal Exchange - Cold Brew"'
$ws.Range("C3").Formula = '="2"'
$ws.Range("D3").Formula = '="71.50"'
$ws.Range("E3").Formula = '="143.00"'

# Convert the formulas above into plain literal text values (keeps default
# style, unlike setting .Value on a numeric-looking string which would
# otherwise coerce the cell to a Number / trigger a text-quote style).
$ws.Range("A2:E3").Copy()
$ws.Range("A2:E3").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# The old rows 4 and 5 no longer exist - shrink the sheet back down.
$ws.Range("A4:E5").Delete()
